$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.424.27'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '1.565.23'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D6").Value = '285.75'
$ws.Range("E6").Value = '  -1.91%  '
$ws.Range("E7").Value = '  -2.72%  '
$ws.Range("D8").Value = '48.55'
$ws.Range("E8").Value = '  -3.09%  '
$ws.Range("D9").Value = '0.3339'
$ws.Range("E9").Value = '  -1.47%  '
$ws.Range("D10").Value = '1.124'
$ws.Range("E10").Value = '  -1.45%  '
$ws.Range("D11").Value = '0.07398'
$ws.Range("E11").Value = '  -2.45%  '
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("D13").Value = '20.73'
$ws.Range("E13").Value = '  -2.83%  '
$ws.Range("E14").Value = '  -0.88%  '
$ws.Range("D15").Value = '6.894'
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("D16").Value = '1.565.80'
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("E17").Value = '  -1.57%  '
$ws.Range("D18").Value = '88.14'
$ws.Range("E18").Value = '  -3.00%  '
$ws.Range("D19").Value = '0.06697'
$ws.Range("E19").Value = '  -0.62%  '
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").Value = '6.348'
$ws.Range("E21").Value = '  +0.97%  '
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("D23").Value = '12.02'
$ws.Range("E23").Value = '  -0.92%  '
$ws.Range("D24").Value = '22.422.82'
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").Value = '2.401'
$ws.Range("E25").Value = '  +2.88%  '
$ws.Range("D26").Value = '2.561'
$ws.Range("E26").Value = '  -3.96%  '
$ws.Range("D27").Value = '150.19'
$ws.Range("E27").Value = '  +1.12%  '
$ws.Range("D28").Value = '19.31'
$ws.Range("E28").Value = '  -3.93%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("D30").Value = '123.60'
$ws.Range("D31").Value = '1.740.15'
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("E32").Value = '  +1.24%  '
$ws.Range("B33").Value = 'WEMIXTOKEN'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D33").Value = '1.999'
$ws.Range("E33").Value = '  +1.18%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '6.099'
$ws.Range("E34").Value = '  -1.03%  '
$ws.Range("D35").Value = '9.792'
$ws.Range("E35").Value = '  -0.41%  '
$ws.Range("E36").Value = '  -1.34%  '
$ws.Range("D37").Value = '0.02406'
$ws.Range("E37").Value = '  -2.75%  '
$ws.Range("D38").Value = '0.2225'
$ws.Range("E38").Value = '  -2.49%  '
$ws.Range("D39").Value = '0.06382'
$ws.Range("E39").Value = '  -2.18%  '
$ws.Range("D40").Value = '1.295'
$ws.Range("E40").Value = '  -6.06%  '
$ws.Range("D41").Value = '5.328'
$ws.Range("E41").Value = '  -2.47%  '
$ws.Range("D42").Value = '11.13'
$ws.Range("E42").Value = '  -1.18%  '
$ws.Range("D43").Value = '0.6084'
$ws.Range("E43").Value = '  -2.25%  '
$ws.Range("D44").Value = '1.001'
$ws.Range("D45").Value = '13.79'
$ws.Range("E45").Value = '  -1.51%  '
$ws.Range("D46").Value = '3.765'
$ws.Range("E46").Value = '  -1.16%  '
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("D48").Value = '2.018'
$ws.Range("E48").Value = '  -2.67%  '
$ws.Range("D49").Value = '123.96'
$ws.Range("E49").Value = '  -4.24%  '
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("D51").Value = '0.07202'
$ws.Range("E51").Value = '  -1.57%  '
